$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 70, shifting existing rows 70-129 down to 71-130.
$ws.Rows("70:70").Insert()

# Populate the newly inserted row 70 with the new data record.
$ws.Range("A70").Value = 6
$ws.Range("B70").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C70").Value = "Metropolitana"
$ws.Range("D70").Value = 45068
$ws.Range("E70").Value = 13
$ws.Range("F70").Value = 100114007
$ws.Range("G70").Value = "Jengibre"
$ws.Range("H70").Value = "Sin especificar"
$ws.Range("I70").Value = "Primera"
$ws.Range("J70").Value = 200
$ws.Range("K70").Value = 13000
$ws.Range("L70").Value = 15000
$ws.Range("M70").Value = 13800
$ws.Range("N70").Value = "$/caja 13 kilos"
$ws.Range("O70").Value = "Perú"
$ws.Range("P70").Value = 1062
$ws.Range("Q70").Value = 13
$ws.Range("R70").Value = "Hortaliza"
